# Update index.html / Excel data: insert a new field report at row 11
# ("LACROZE FEDERICO" pole, case -14) and push the existing rows 11-21
# down by one (row 22's old content is discarded, matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11-21 down into rows 12-22 (process bottom-up so a source row
# is never overwritten before it's copied from).
for ($r = 21; $r -ge 11; $r--) {
    $src = $ws.Range("A" + $r + ":N" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":N" + ($r + 1))
    $src.Copy($dst)
}

# Write the brand-new case into row 11.
$ws.Range("A11:E11").NumberFormat = "@"
$ws.Range("H11:L11").NumberFormat = "@"

$ws.Range("A11").Value = "-14"
$ws.Range("B11").Value = "1/23/2024"
$ws.Range("C11").Value = "LACROZE FEDERICO ,AV. /ALT/ 3057"
$ws.Range("D11").Value = "106581 - COLEGIALES"
$ws.Range("E11").Value = "778723976"
$ws.Range("F11").Value = "GESTION TELECENTRO"
$ws.Range("G11").Value = "Pendiente"
$ws.Range("H11").Value = "Fede. Lacroze 3057 cambiar columna 114 efectuar transferencias +fuente y cdo. Base corroida"
$ws.Range("I11").Value = "1"
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = -58.44873
$ws.Range("N11").Value = -34.576501
